$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 807.8333
$ws.Range("I2").Value = 729.4
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 729.4
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -616.4
$ws.Range("N2").Value = -1426

$ws.Range("H40").Value = 5366.727
$ws.Range("I40").Value = 3405.25
$ws.Range("J40").Value = 10597.333
$ws.Range("K40").Value = 3405.25
$ws.Range("L40").Value = 10597.333
$ws.Range("M40").Value = -3230.25
$ws.Range("N40").Value = -10947.333

$ws.Range("H41").Value = 1219.3529
$ws.Range("J41").Value = 1016.375
$ws.Range("L41").Value = 1016.375
$ws.Range("N41").Value = -1896.375

$ws.Range("H48").Value = 19000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H51").Value = 2792.8572
$ws.Range("I51").Value = 1275
$ws.Range("J51").Value = 3400
$ws.Range("K51").Value = 1275
$ws.Range("L51").Value = 3400
$ws.Range("M51").Value = -791
$ws.Range("N51").Value = -4368

$ws.Range("H56").Value = 19000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H111").Value = 1150
$ws.Range("I111").Value = 950
$ws.Range("J111").Value = 1250
$ws.Range("K111").Value = 2850
$ws.Range("L111").Value = 3750
$ws.Range("M111").Value = 217
$ws.Range("N111").Value = -9884

$ws.Range("H113").Value = 3607.1
$ws.Range("I113").Value = 2454.6
$ws.Range("K113").Value = 2454.6
$ws.Range("M113").Value = 799.4000000000001

$ws.Range("H135").Value = 2886.6667
$ws.Range("I135").Value = 2886.6667
$ws.Range("K135").Value = 25980.0003
$ws.Range("M135").Value = -23445.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 144.5
$ws.Range("I5").Value = 140.23077
$ws.Range("K5").Value = 140.23077
$ws.Range("M5").Value = -28.23077000000001

$ws.Range("H32").Value = 5544.1616
$ws.Range("I32").Value = 1573.9215
$ws.Range("J32").Value = 17454.883
$ws.Range("K32").Value = 1573.9215
$ws.Range("L32").Value = 17454.883
$ws.Range("M32").Value = -1286.9215
$ws.Range("N32").Value = -18028.883

$ws.Range("H76").Value = 500075000
$ws.Range("J76").Value = 500075000
$ws.Range("L76").Value = 500075000
$ws.Range("N76").Value = -500075676

$ws.Range("H79").Value = 500075000
$ws.Range("J79").Value = 500075000
$ws.Range("L79").Value = 500075000
$ws.Range("N79").Value = -500077340

$ws.Range("H110").Value = 1623.4546
$ws.Range("I110").Value = 881.8
$ws.Range("J110").Value = 2241.5
$ws.Range("K110").Value = 881.8
$ws.Range("L110").Value = 2241.5
$ws.Range("M110").Value = 1163.2
$ws.Range("N110").Value = -6331.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 144.5
$ws.Range("I4").Value = 140.23077
$ws.Range("K4").Value = 140.23077
$ws.Range("M4").Value = -25.23077000000001

$ws.Range("H22").Value = 62631.47
$ws.Range("I22").Value = 95014
$ws.Range("K22").Value = 95014
$ws.Range("M22").Value = -94841

$ws.Range("H29").Value = 658
$ws.Range("I29").Value = 316
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 316
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = -27
$ws.Range("N29").Value = -1578

$ws.Range("H75").Value = 8585.6
$ws.Range("I75").Value = 8585.6
$ws.Range("K75").Value = 8585.6
$ws.Range("M75").Value = -7649.6

$ws.Range("H78").Value = 8585.6
$ws.Range("I78").Value = 8585.6
$ws.Range("K78").Value = 25756.8
$ws.Range("M78").Value = -21076.8

$ws.Range("H134").Value = 3546.4412
$ws.Range("J134").Value = 6748.6665
$ws.Range("L134").Value = 20245.9995
$ws.Range("N134").Value = -25315.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 336.875
$ws.Range("I22").Value = 219
$ws.Range("J22").Value = 533.3333
$ws.Range("K22").Value = 219
$ws.Range("L22").Value = 533.3333
$ws.Range("M22").Value = 131
$ws.Range("N22").Value = -1233.3333

$ws.Range("H31").Value = 3091.1428
$ws.Range("I31").Value = 2394.5454
$ws.Range("K31").Value = 2394.5454
$ws.Range("M31").Value = -2099.5454

$ws.Range("H34").Value = 3091.1428
$ws.Range("I34").Value = 2394.5454
$ws.Range("K34").Value = 2394.5454
$ws.Range("M34").Value = -2192.5454

$ws.Range("H62").Value = 2486.125
$ws.Range("J62").Value = 2132.6667
$ws.Range("L62").Value = 2132.6667
$ws.Range("N62").Value = -3380.6667

$ws.Range("H65").Value = 2486.125
$ws.Range("J65").Value = 2132.6667
$ws.Range("L65").Value = 10663.3335
$ws.Range("N65").Value = -16903.3335

$ws.Range("H99").Value = 10171287
$ws.Range("I99").Value = 15875510
$ws.Range("K99").Value = 15875510
$ws.Range("M99").Value = -15874012

$ws.Range("H126").Value = 10171287
$ws.Range("I126").Value = 15875510
$ws.Range("K126").Value = 47626530
$ws.Range("M126").Value = -47624060

$ws.Range("H132").Value = 1828611
$ws.Range("J132").Value = 1001844.06
$ws.Range("L132").Value = 3005532.18
$ws.Range("N132").Value = -3010592.18

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 897.3
$ws.Range("I5").Value = 797.4545000000001
$ws.Range("J5").Value = 1019.3333
$ws.Range("K5").Value = 2392.3635
$ws.Range("L5").Value = 3057.9999
$ws.Range("M5").Value = -2280.3635
$ws.Range("N5").Value = -3281.9999

$ws.Range("H113").Value = 2431916
$ws.Range("I113").Value = 1220.8
$ws.Range("J113").Value = 3242147.8
$ws.Range("K113").Value = 3662.4
$ws.Range("L113").Value = 9726443.399999999
$ws.Range("M113").Value = -1492.4
$ws.Range("N113").Value = -9730783.399999999

$ws.Range("H129").Value = 988.6667
$ws.Range("J129").Value = 1631.6666
$ws.Range("L129").Value = 4894.9998
$ws.Range("N129").Value = -14894.9998

$ws.Range("H135").Value = 897.3
$ws.Range("I135").Value = 797.4545000000001
$ws.Range("J135").Value = 1019.3333
$ws.Range("K135").Value = 7177.0905
$ws.Range("L135").Value = 9173.9997
$ws.Range("M135").Value = -4642.0905
$ws.Range("N135").Value = -14243.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 14443.5
$ws.Range("J96").Value = 14443.5
$ws.Range("L96").Value = 14443.5
$ws.Range("N96").Value = -19935.5

$ws.Range("H98").Value = 8250
$ws.Range("J98").Value = 8250
$ws.Range("L98").Value = 8250
$ws.Range("N98").Value = -14240

$ws.Range("H126").Value = 4244.7
$ws.Range("I126").Value = 3179.4
$ws.Range("K126").Value = 9538.200000000001
$ws.Range("M126").Value = -7068.200000000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1217
$ws.Range("I22").Value = 1373.3846
$ws.Range("J22").Value = 962.875
$ws.Range("K22").Value = 1373.3846
$ws.Range("L22").Value = 962.875
$ws.Range("M22").Value = -1078.3846
$ws.Range("N22").Value = -1552.875

$ws.Range("H24").Value = 2925
$ws.Range("J24").Value = 2925
$ws.Range("L24").Value = 2925
$ws.Range("N24").Value = -3611

$ws.Range("H27").Value = 1217
$ws.Range("I27").Value = 1373.3846
$ws.Range("J27").Value = 962.875
$ws.Range("K27").Value = 1373.3846
$ws.Range("L27").Value = 962.875
$ws.Range("M27").Value = -1266.3846
$ws.Range("N27").Value = -1176.875

$ws.Range("H30").Value = 4994
$ws.Range("J30").Value = 4994
$ws.Range("L30").Value = 4994
$ws.Range("N30").Value = -5210

$ws.Range("H132").Value = 2776.375
$ws.Range("J132").Value = 2005
$ws.Range("L132").Value = 6015
$ws.Range("N132").Value = -11075

$ws.Range("H136").Value = 2869.7144
$ws.Range("I136").Value = 2366.5
$ws.Range("K136").Value = 7099.5
$ws.Range("M136").Value = -4549.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 16500
$ws.Range("J20").Value = 16500
$ws.Range("L20").Value = 16500
$ws.Range("N20").Value = -16980

$ws.Range("H87").Value = 76499.5
$ws.Range("J87").Value = 76499.5
$ws.Range("L87").Value = 76499.5
$ws.Range("N87").Value = -78995.5

$ws.Range("H90").Value = 76499.5
$ws.Range("J90").Value = 76499.5
$ws.Range("L90").Value = 229498.5
$ws.Range("N90").Value = -241978.5

$ws.Range("H114").Value = 10000
$ws.Range("J114").Value = 10000
$ws.Range("L114").Value = 10000
$ws.Range("N114").Value = -18678

$ws.Range("H115").Value = 66648
$ws.Range("J115").Value = 66648
$ws.Range("L115").Value = 66648
$ws.Range("N115").Value = -69782

$ws.Range("H122").Value = 1408.0769
$ws.Range("I122").Value = 1140.6
$ws.Range("J122").Value = 2299.6667
$ws.Range("K122").Value = 3421.8
$ws.Range("L122").Value = 6899.000100000001
$ws.Range("M122").Value = -971.7999999999997
$ws.Range("N122").Value = -11799.0001

$ws.Range("H126").Value = 2975.8
$ws.Range("I126").Value = 3001.3333
$ws.Range("K126").Value = 9003.999899999999
$ws.Range("M126").Value = -6533.999899999999

$ws.Range("H135").Value = 148447.33
$ws.Range("J135").Value = 148447.33
$ws.Range("L135").Value = 148447.33
$ws.Range("N135").Value = -158587.33
